{"js": "// Update the date line and the 25 division-problem table cells to the\n// new set of values, per the commit's regenerated worksheet output.\nconst replacements = [\n  [\"2024-07-31 Wednesday\", \"2024-08-01 Thursday\"],\n  [\"53\u00f76=8, 5\", \"75\u00f79=8, 3\"],\n  [\"82\u00f74=20, 2\", \"75\u00f76=12, 3\"],\n  [\"74\u00f77=10, 4\", \"28\u00f79=3, 1\"],\n  [\"79\u00f74=19, 3\", \"75\u00f77=10, 5\"],\n  [\"16\u00f74=4, 0\", \"66\u00f74=16, 2\"],\n  [\"33\u00f72=16, 1\", \"97\u00f75=19, 2\"],\n  [\"31\u00f76=5, 1\", \"49\u00f75=9, 4\"],\n  [\"26\u00f77=3, 5\", \"30\u00f79=3, 3\"],\n  [\"25\u00f73=8, 1\", \"84\u00f76=14, 0\"],\n  [\"40\u00f74=10, 0\", \"74\u00f74=18, 2\"],\n  [\"61\u00f72=30, 1\", \"35\u00f74=8, 3\"],\n  [\"79\u00f78=9, 7\", \"14\u00f79=1, 5\"],\n  [\"67\u00f75=13, 2\", \"43\u00f79=4, 7\"],\n  [\"72\u00f75=14, 2\", \"51\u00f73=17, 0\"],\n  [\"57\u00f78=7, 1\", \"69\u00f72=34, 1\"],\n  [\"92\u00f74=23, 0\", \"25\u00f79=2, 7\"],\n  [\"48\u00f76=8, 0\", \"67\u00f75=13, 2\"],\n  [\"17\u00f73=5, 2\", \"95\u00f78=11, 7\"],\n  [\"61\u00f75=12, 1\", \"77\u00f75=15, 2\"],\n  [\"83\u00f73=27, 2\", \"44\u00f78=5, 4\"],\n  [\"22\u00f74=5, 2\", \"99\u00f77=14, 1\"],\n  [\"28\u00f74=7, 0\", \"75\u00f72=37, 1\"],\n  [\"93\u00f78=11, 5\", \"24\u00f74=6, 0\"],\n  [\"30\u00f78=3, 6\", \"35\u00f72=17, 1\"],\n  [\"83\u00f74=20, 3\", \"32\u00f77=4, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem table cells to the\n# new set of values, per the commit's regenerated worksheet output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-31 Wednesday\", \"2024-08-01 Thursday\"),\n    @(\"53\u00f76=8, 5\", \"75\u00f79=8, 3\"),\n    @(\"82\u00f74=20, 2\", \"75\u00f76=12, 3\"),\n    @(\"74\u00f77=10, 4\", \"28\u00f79=3, 1\"),\n    @(\"79\u00f74=19, 3\", \"75\u00f77=10, 5\"),\n    @(\"16\u00f74=4, 0\", \"66\u00f74=16, 2\"),\n    @(\"33\u00f72=16, 1\", \"97\u00f75=19, 2\"),\n    @(\"31\u00f76=5, 1\", \"49\u00f75=9, 4\"),\n    @(\"26\u00f77=3, 5\", \"30\u00f79=3, 3\"),\n    @(\"25\u00f73=8, 1\", \"84\u00f76=14, 0\"),\n    @(\"40\u00f74=10, 0\", \"74\u00f74=18, 2\"),\n    @(\"61\u00f72=30, 1\", \"35\u00f74=8, 3\"),\n    @(\"79\u00f78=9, 7\", \"14\u00f79=1, 5\"),\n    @(\"67\u00f75=13, 2\", \"43\u00f79=4, 7\"),\n    @(\"72\u00f75=14, 2\", \"51\u00f73=17, 0\"),\n    @(\"57\u00f78=7, 1\", \"69\u00f72=34, 1\"),\n    @(\"92\u00f74=23, 0\", \"25\u00f79=2, 7\"),\n    @(\"48\u00f76=8, 0\", \"67\u00f75=13, 2\"),\n    @(\"17\u00f73=5, 2\", \"95\u00f78=11, 7\"),\n    @(\"61\u00f75=12, 1\", \"77\u00f75=15, 2\"),\n    @(\"83\u00f73=27, 2\", \"44\u00f78=5, 4\"),\n    @(\"22\u00f74=5, 2\", \"99\u00f77=14, 1\"),\n    @(\"28\u00f74=7, 0\", \"75\u00f72=37, 1\"),\n    @(\"93\u00f78=11, 5\", \"24\u00f74=6, 0\"),\n    @(\"30\u00f78=3, 6\", \"35\u00f72=17, 1\"),\n    @(\"83\u00f74=20, 3\", \"32\u00f77=4, 4\")\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $after, 2) | Out-Null\n}\n"}
